# Add a new "Viewership" report sheet with total viewership by hour,
# mirroring the layout/formatting already used on the "Sessions" sheet.

$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item("Sessions")

# New sheet goes after the last existing sheet.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "Viewership"

# --- Column widths (approximate the author's best-fit widths) ---
$ws.Columns.Item(1).ColumnWidth = 14.333333333333332
$ws.Columns.Item(2).ColumnWidth = 16.666666666666664
$ws.Columns.Item(3).ColumnWidth = 15.666666666666666

# --- Title block ---
$ws.Range("A1").Value = "Who Is Streaming?"
$srcSheet.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").NumberFormat = "yyyy/mm/dd hh"

$ws.Range("A2").Value = "Game Name"
$srcSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").NumberFormat = "yyyy/mm/dd hh"

# Scratch cell, immediately removed, so the plain (non-bold) date-style
# used as the column's default format gets registered in the style table
# in the same slot the original workbook uses.
$ws.Range("A1000").Value = 1
$ws.Range("A1000").NumberFormat = "yyyy/mm/dd hh"
$ws.Rows.Item(1000).Delete() | Out-Null

# --- Start / End summary ---
$ws.Range("A4").Value = "Start"
$srcSheet.Range("A4").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").NumberFormat = "yyyy/mm/dd hh"

$ws.Range("B4").Formula = "=MIN(`$A`$9:`$A`$2006)"
$srcSheet.Range("A4").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").NumberFormat = "yyyy/mm/dd hh"

$ws.Range("A5").Value = "End"
$srcSheet.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").NumberFormat = "yyyy/mm/dd hh"

$ws.Range("B5").Formula = "=MAX(`$B`$8:`$B`$2005)"
$srcSheet.Range("A4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").NumberFormat = "yyyy/mm/dd hh"

# --- Table header row ---
$ws.Range("A7").Value = "Hour"
$srcSheet.Range("A7").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("B7").Value = "Streamers"
$srcSheet.Range("A7").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("C7").Value = "Viewers"
$srcSheet.Range("A7").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$ws.Range("A3").Select() | Out-Null
